$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''37.802.03'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '''2.085.67'
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''235.12'
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").Value = '''0.625'
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").Value = '''59.81'
$ws.Range("E7").Value = '  +4.10%  '
$ws.Range("D9").Value = '''0.391'
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("D10").Value = '''0.0792'
$ws.Range("E10").Value = '  +2.26%  '
$ws.Range("E11").Value = '  +2.95%  '
$ws.Range("D12").Value = '''2.390.50'
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").Value = '''14.70'
$ws.Range("E13").Value = '  +1.48%  '
$ws.Range("D14").Value = '''21.46'
$ws.Range("E14").Value = '  +4.13%  '
$ws.Range("D15").Value = '''0.773'
$ws.Range("E15").Value = '  -0.88%  '
$ws.Range("D16").Value = '''5.33'
$ws.Range("E16").Value = '  +2.74%  '
$ws.Range("D17").Value = '''2.080.34'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").Value = '''37.718.15'
$ws.Range("E18").Value = '  +0.91%  '
$ws.Range("E19").Value = '  -3.12%  '
$ws.Range("D20").Value = '''71.65'
$ws.Range("E20").Value = '  +2.59%  '
$ws.Range("D21").Value = '''0.0₃0830'
$ws.Range("E21").Value = '  +1.33%  '
$ws.Range("D22").Value = '''229.03'
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").Value = '''2.42'
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("D26").Value = '''170.51'
$ws.Range("E26").Value = '  +2.16%  '
$ws.Range("D27").Value = '''0.141'
$ws.Range("E27").Value = '  +10.03%  '
$ws.Range("D28").Value = '''9.07'
$ws.Range("E28").Value = '  +2.44%  '
$ws.Range("D29").Value = '''1.44'
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("D30").Value = '''19.54'
$ws.Range("E30").Value = '  +1.98%  '
$ws.Range("E31").Value = '  +1.34%  '
$ws.Range("D32").Value = '''4.74'
$ws.Range("E32").Value = '  +4.19%  '
$ws.Range("D33").Value = '''0.0632'
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("D34").Value = '''4.71'
$ws.Range("E34").Value = '  +2.69%  '
$ws.Range("D35").Value = '''2.53'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").Value = '''3.55'
$ws.Range("E36").Value = '  +6.76%  '
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("D39").Value = '''5.45'
$ws.Range("E39").Value = '  -4.69%  '
$ws.Range("D40").Value = '''0.0989'
$ws.Range("E40").Value = '  +2.53%  '
$ws.Range("D41").Value = '''100.06'
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("D44").Value = '''1.464.99'
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("D45").Value = '''1.18'
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").Value = '''4.22'
$ws.Range("E46").Value = '  +2.95%  '
$ws.Range("D47").Value = '''16.20'
$ws.Range("E47").Value = '  +5.53%  '
$ws.Range("D48").Value = '''1.07'
$ws.Range("E48").Value = '  +3.44%  '
$ws.Range("D49").Value = '''7.48'
$ws.Range("E49").Value = '  +3.04%  '
$ws.Range("E50").Value = '  +2.64%  '
$ws.Range("D51").Value = '''47.50'
$ws.Range("E51").Value = '  +6.11%  '
